$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "settings": rename the default bot name, and add a new
# CHATWORK_API_TOKEN setting row (inserted after SLACK_ICON_EMOJI, before
# TIME_INTERVAL).
# ---------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("settings")

$wsSettings.Range("B2").Value2 = "EndlessReminder"

$wsSettings.Rows("5").Copy()
$wsSettings.Rows("6").Insert()
$wsSettings.Range("A6").Value2 = "CHATWORK_API_TOKEN"
$wsSettings.Range("B6").Value2 = ""

$wsSettings.Activate()
[void]$wsSettings.Range("A10").Select()

# ---------------------------------------------------------------------
# Sheet "main": add a new scheduling row (row 3) that duplicates row 2's
# style/content but schedules for Friday ("Fri") instead of every day,
# and is not currently completed (not-done checkbox = FALSE).
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("main")

$wsMain.Rows("2").Copy()
$wsMain.Rows("3").Insert()
$wsMain.Rows("3").RowHeight = 56

$wsMain.Range("A3").Value2 = 2
$wsMain.Range("D3").Value2 = "Fri"
$wsMain.Range("E3").Value2 = $false

$wsMain.Activate()
[void]$wsMain.Range("A3").Select()
